# menambahkan fungsi bulk toko csv
#
# Adds a new "Status" column to Sheet1 (inserted as column C, pushing the
# existing Alamat / Link GMAP / latitude / longitude / Icon columns one to
# the right) and fills it in with each outlet's status, as produced by the
# bulk "toko csv" import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column in front of the current column C ("Alamat").
# Everything from C onward (Alamat, Link GMAP, latitude, longitude, Icon)
# shifts one column to the right automatically.
$ws.Columns("C:C").Insert()

# New header cell - plain/unformatted text like the bulk-imported data,
# not bold like the rest of the header row.
$ws.Range("C1").Style = "Normal"
$ws.Range("C1").Value = "Status"

# Per-row outlet status values (row 1 is the header row).
$status = [ordered]@{
    2  = "Active"
    3  = "Active"
    4  = "Active"
    5  = "Active"
    6  = "Active"
    7  = "Active"
    8  = "Active"
    9  = "Active"
    10 = "Active"
    11 = "Active"
    12 = "Active"
    13 = "Active"
    14 = "Active"
    15 = "Active"
    16 = "Active"
    17 = "Active"
    18 = "Active"
    19 = "Active"
    20 = "Active"
    21 = "Active"
    22 = "Active"
    23 = "on-Survey"
    24 = "Active"
    25 = "Active"
    26 = "on-Survey"
    27 = "Active"
    28 = "Active"
    29 = "on-Survey"
    30 = "on-Survey"
}

foreach ($row in $status.Keys) {
    $ws.Cells.Item($row, 3).Value = $status[$row]
}

# Restore a sensible active selection, same spirit as the source edit.
[void]$ws.Range("E36").Select()
